$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1D NEW")
$ws.Activate()

# --- Add three new rows to the Table43 table (expands table range B3:M33 -> B3:M36) ---
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Row 34: 1dmockanderrors31.csv ---
$ws.Range("B34").Value = "1dmockanderrors31.csv"
$ws.Range("C34").Value = 50
$ws.Range("D34").Value = 1000
$ws.Range("E34").Value = 0.3
$ws.Range("F34").Value = 0.075
$ws.Range("G34").Value = 360
$ws.Range("H34").Value = 1
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 5
$ws.Range("L34").Value = 1

# --- Row 35: 1dmockanderrors32.csv ---
$ws.Range("B35").Value = "1dmockanderrors32.csv"
$ws.Range("C35").Value = 50
$ws.Range("D35").Value = 1000
$ws.Range("E35").Value = 0.3
$ws.Range("F35").Value = 0.075
$ws.Range("G35").Value = 360
$ws.Range("H35").Value = 1
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 1
$ws.Range("M35").Value = "noiseless"

# --- Row 36: 1dmockanderrors33.csv ---
$ws.Range("B36").Value = "1dmockanderrors33.csv"
$ws.Range("C36").Value = 50
$ws.Range("D36").Value = 1000
$ws.Range("E36").Value = "[Two equal peaks at 0.182 and 0.273]"
$ws.Range("G36").Value = 366
$ws.Range("H36").Value = 1
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 1
$ws.Range("M36").Value = "realistic but larger pixel pitch"
$ws.Range("F36").Value = "[modulated by a gaussian with FWHM = 30mm]"

# --- Update sheet view: scroll down + select F37 (new last comment cell) ---
$ws.Range("A16").Select() | Out-Null
$ws.Range("F37").Select() | Out-Null
